# Update "想去人数" (column F) figures across the sheets, per the upstream
# gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 8119
$ws.Range("F3").Value = 121
$ws.Range("F5").Value = 31035
$ws.Range("F8").Value = 704
$ws.Range("F12").Value = 805
$ws.Range("F13").Value = 54
$ws.Range("F14").Value = 608
$ws.Range("F15").Value = 397
$ws.Range("F17").Value = 559
$ws.Range("F18").Value = 150
$ws.Range("F19").Value = 415
$ws.Range("F23").Value = 710
$ws.Range("F24").Value = 2350
$ws.Range("F25").Value = 838
$ws.Range("F26").Value = 69
$ws.Range("F27").Value = 1090
$ws.Range("F29").Value = 636
$ws.Range("F30").Value = 1078

# --- 演出 (sheet 2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 296
$ws.Range("F5").Value = 317
$ws.Range("F8").Value = 31

# --- 本地生活 (sheet 3) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 527

# --- 全部类型 (sheet 4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 527
$ws.Range("F3").Value = 8119
$ws.Range("F4").Value = 121
$ws.Range("F6").Value = 296
$ws.Range("F7").Value = 31035
$ws.Range("F10").Value = 704
$ws.Range("F16").Value = 317
$ws.Range("F18").Value = 805
$ws.Range("F19").Value = 54
$ws.Range("F20").Value = 608
$ws.Range("F21").Value = 397
$ws.Range("F24").Value = 31
$ws.Range("F27").Value = 559
$ws.Range("F28").Value = 150
$ws.Range("F29").Value = 415
$ws.Range("F33").Value = 710
$ws.Range("F34").Value = 2350
$ws.Range("F35").Value = 838
$ws.Range("F36").Value = 69
$ws.Range("F40").Value = 636
$ws.Range("F41").Value = 1078
